$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.311.97'
$ws.Range('E2').Value = '  +0.29%  '
$ws.Range('D3').Value = '3.344.12'
$ws.Range('E3').Value = '  +0.65%  '
$ws.Range('D4').Value = "'0.998"
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').Value = "'585.30"
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +4.96%  '
$ws.Range('D6').Value = "'185.74"
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -1.40%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('E8').Value = '  -1.34%  '
$ws.Range('E9').Value = '  -0.92%  '
$ws.Range('D10').Value = "'0.582"
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -0.67%  '
$ws.Range('D11').Value = "'46.98"
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -1.07%  '
$ws.Range('E12').Value = '  -0.28%  '
$ws.Range('D13').Value = "'654.58"
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +7.45%  '
$ws.Range('E14').Value = '  -2.47%  '
$ws.Range('D15').Value = '3.627.96'
$ws.Range('E15').Value = '  -5.78%  '
$ws.Range('D16').Value = '66.406.71'
$ws.Range('E16').Value = '  +0.45%  '
$ws.Range('E17').Value = '  -0.21%  '
$ws.Range('E18').Value = '  -0.80%  '
$ws.Range('D19').Value = '3.339.25'
$ws.Range('E19').Value = '  +0.98%  '
$ws.Range('E20').Value = '  +0.15%  '
$ws.Range('E21').Value = '  -1.24%  '
$ws.Range('D22').Value = "'17.69"
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -4.16%  '
$ws.Range('E23').Value = '  -0.41%  '
$ws.Range('D24').Value = "'100.30"
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -0.23%  '
$ws.Range('E25').Value = '  +0.74%  '
$ws.Range('E26').Value = '  +0.70%  '
$ws.Range('D27').Value = "'9.61"
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -0.07%  '
$ws.Range('D28').Value = "'32.09"
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +5.45%  '
$ws.Range('D29').Value = "'8.54"
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -2.16%  '
$ws.Range('D30').Value = "'6.85"
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +1.25%  '
$ws.Range('D31').Value = "'601.23"
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +3.14%  '
$ws.Range('E32').Value = '  -0.12%  '
$ws.Range('E33').Value = '  -0.06%  '
$ws.Range('D34').Value = '3.878.70'
$ws.Range('E34').Value = '  +4.35%  '
$ws.Range('E35').Value = '  +0.37%  '
$ws.Range('E36').Value = '  -0.01%  '
$ws.Range('D37').Value = "'56.45"
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -1.12%  '
$ws.Range('D38').Value = "'2.74"
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +1.92%  '
$ws.Range('E39').Value = '  -0.66%  '
$ws.Range('D40').Value = "'33.11"
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -2.75%  '
$ws.Range('E41').Value = '  -2.42%  '
$ws.Range('E42').Value = '  -2.94%  '
$ws.Range('E43').Value = '  -0.20%  '
$ws.Range('E44').Value = '  -0.33%  '
$ws.Range('E45').Value = '  -1.18%  '
$ws.Range('E46').Value = '  -1.56%  '
$ws.Range('B47').Value = 'FirstDigitalUSD'
$ws.Range('C47').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D47').Value = "'1.00"
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +0.42%  '
$ws.Range('B48').Value = 'ThetaToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D48').Value = "'2.56"
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -1.81%  '
$ws.Range('D49').Value = "'2.87"
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -16.17%  '
$ws.Range('E50').Value = '  +6.62%  '
$ws.Range('D51').Value = "'130.22"
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +5.32%  '
